# Update cryptos price list cells (Coin, Link, Price, Volume(1h)) for Fri Jun 30 2023 refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.693.79"
$ws.Range("E2").Value = "  +1.25%  "
$ws.Range("D3").Value = "'1.878.18"
$ws.Range("E3").Value = "  +1.97%  "
$ws.Range("D4").Value = "'0.9996"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'237.25"
$ws.Range("E5").Value = "  +1.85%  "
$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.4739"
$ws.Range("E7").Value = "  +2.09%  "
$ws.Range("D8").Value = "'0.2825"
$ws.Range("E8").Value = "  +3.58%  "
$ws.Range("D9").Value = "'0.06504"
$ws.Range("E9").Value = "  +3.77%  "
$ws.Range("D10").Value = "'18.63"
$ws.Range("E10").Value = "  +14.69%  "
$ws.Range("D11").Value = "'1.878.77"
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").Value = "'0.07571"
$ws.Range("E12").Value = "  +2.03%  "
$ws.Range("D13").Value = "'95.56"
$ws.Range("E13").Value = "  +14.26%  "
$ws.Range("D14").Value = "'5.080"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").Value = "'0.6477"
$ws.Range("E15").Value = "  +4.47%  "
$ws.Range("D16").Value = "'306.51"
$ws.Range("E16").Value = "  +34.51%  "
$ws.Range("D17").Value = "'30.696.31"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.07"
$ws.Range("E18").Value = "  +5.93%  "
$ws.Range("B19").Value = "Dai"
$ws.Range("C19").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D19").Value = "'0.9988"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'0.000007533"
$ws.Range("E20").Value = "  +3.50%  "
$ws.Range("D21").Value = "'2.121.82"
$ws.Range("E21").Value = "  +2.45%  "
$ws.Range("D22").Value = "'0.9992"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "'5.142"
$ws.Range("E23").Value = "  +4.80%  "
$ws.Range("D24").Value = "'6.156"
$ws.Range("E24").Value = "  +5.06%  "
$ws.Range("D25").Value = "'169.01"
$ws.Range("E25").Value = "  +2.85%  "
$ws.Range("D26").Value = "'9.238"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("D27").Value = "'19.73"
$ws.Range("E27").Value = "  +10.99%  "
$ws.Range("D28").Value = "'1.947"
$ws.Range("E28").Value = "  +4.25%  "
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").Value = "'1.352"
$ws.Range("E30").Value = "  -1.26%  "
$ws.Range("D31").Value = "'4.161"
$ws.Range("E31").Value = "  +2.11%  "
$ws.Range("D32").Value = "'3.941"
$ws.Range("E32").Value = "  +3.48%  "
$ws.Range("D33").Value = "'0.05034"
$ws.Range("E33").Value = "  +3.97%  "
$ws.Range("D34").Value = "'1.172"
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("D35").Value = "'0.7193"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("D36").Value = "'2.706"
$ws.Range("E36").Value = "  +0.46%  "
$ws.Range("D37").Value = "'0.01916"
$ws.Range("E37").Value = "  +2.05%  "
$ws.Range("D38").Value = "'2.699"
$ws.Range("E38").Value = "  +1.76%  "
$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").Value = "'0.8970"
$ws.Range("E39").Value = "  +1.08%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").Value = "'2.042"
$ws.Range("E40").Value = "  +6.29%  "
$ws.Range("D41").Value = "'106.97"
$ws.Range("E41").Value = "  +1.93%  "
$ws.Range("D42").Value = "'0.9986"
$ws.Range("E42").Value = "  -0.29%  "
$ws.Range("D43").Value = "'0.4187"
$ws.Range("E43").Value = "  +4.35%  "
$ws.Range("D44").Value = "'5.580"
$ws.Range("E44").Value = "  +0.65%  "
$ws.Range("D45").Value = "'7.315"
$ws.Range("E45").Value = "  +3.66%  "
$ws.Range("D46").Value = "'65.07"
$ws.Range("E46").Value = "  +7.60%  "
$ws.Range("D47").Value = "'8.950"
$ws.Range("E47").Value = "  +4.09%  "
$ws.Range("D48").Value = "'0.1215"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Value = "'34.52"
$ws.Range("E49").Value = "  +4.15%  "
$ws.Range("D50").Value = "'0.05579"
$ws.Range("E50").Value = "  +1.40%  "
$ws.Range("E51").Value = "  +2.28%  "
